$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.02
$ws.Range("U3").Value = 1.04
$ws.Range("J4").Value = 3.45
$ws.Range("Q4").Value = 1.84
$ws.Range("R4").Value = 1.31
$ws.Range("F5").Value = 2.02
$ws.Range("H5").Value = 3.6
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 1.01
$ws.Range("T5").Value = 1.54
$ws.Range("R6").Value = 1.79
$ws.Range("S6").Value = 2
$ws.Range("R7").Value = 1.18
$ws.Range("T7").Value = 1.76
$ws.Range("P9").Value = 1.72
$ws.Range("Q9").Value = 2.32
$ws.Range("T9").Value = 1.99
$ws.Range("AB9").Value = 8.199999999999999
$ws.Range("AL9").Value = 48
$ws.Range("AM9").Value = 140
$ws.Range("AN9").Value = 23
$ws.Range("L10").Value = 1.34
$ws.Range("P11").Value = 2.66
$ws.Range("Q11").Value = 1.59
$ws.Range("R11").Value = 1.66
$ws.Range("T11").Value = 1.67
$ws.Range("AA11").Value = 150
$ws.Range("AC11").Value = 11
$ws.Range("F12").Value = 3.4
$ws.Range("I12").Value = 2.26
$ws.Range("R12").Value = 1.47
$ws.Range("S12").Value = 3
$ws.Range("T12").Value = 1.67
$ws.Range("V12").Value = 1.79
$ws.Range("F14").Value = 2.48
$ws.Range("G14").Value = 2.84
$ws.Range("H14").Value = 2.86
$ws.Range("I14").Value = 3.3
$ws.Range("J14").Value = 3.15
$ws.Range("K14").Value = 3.65
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 3.25
$ws.Range("O14").Value = 1.35
$ws.Range("P14").Value = 1.78
$ws.Range("Q14").Value = 2.04
$ws.Range("R14").Value = 1.27
$ws.Range("S14").Value = 3.3
$ws.Range("T14").Value = 1.79
$ws.Range("U14").Value = 2.04
$ws.Range("V14").Value = 1.43
$ws.Range("W14").Value = 1.55
$ws.Range("Y14").Value = 16
$ws.Range("Z14").Value = 30
$ws.Range("AB14").Value = 15
$ws.Range("AC14").Value = 11
$ws.Range("AD14").Value = 20
$ws.Range("AF14").Value = 25
$ws.Range("AG14").Value = 18
$ws.Range("AH14").Value = 26
$ws.Range("AK14").Value = 46
$ws.Range("AN14").Value = 40
$ws.Range("AO14").Value = 50
$ws.Range("K15").Value = 3.65
$ws.Range("G19").Value = 1.52
$ws.Range("H19").Value = 6.8
$ws.Range("J22").Value = 3.6
$ws.Range("Q22").Value = 1.86
$ws.Range("H23").Value = 5.8
$ws.Range("P24").Value = 2.48
$ws.Range("F25").Value = 2.5
$ws.Range("R25").Value = 1.09
$ws.Range("G26").Value = 1.69
$ws.Range("J26").Value = 4
$ws.Range("W26").Value = 2.44
$ws.Range("F30").Value = 3.3
$ws.Range("K30").Value = 5.3
$ws.Range("AJ30").Value = 980
$ws.Range("F31").Value = 1.46
$ws.Range("G31").Value = 1.47
$ws.Range("I31").Value = 8.800000000000001
$ws.Range("N31").Value = 3.95
$ws.Range("R31").Value = 1.38
$ws.Range("S31").Value = 3.4
$ws.Range("U31").Value = 1.79
$ws.Range("V31").Value = 1.12
$ws.Range("Z31").Value = 75
$ws.Range("AA31").Value = 390
$ws.Range("AB31").Value = 7.6
$ws.Range("AK31").Value = 16
$ws.Range("AL31").Value = 46
$ws.Range("AM31").Value = 180
$ws.Range("AO31").Value = 190
$ws.Range("AA32").Value = 340
$ws.Range("AN32").Value = 5.1
$ws.Range("Q33").Value = 1.96
